$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dimension/measure labels that describe each column of the data dictionary.
# These columns were "dimension" typed and are now re-typed as "measure", and a couple
# of refArea columns now carry the curated field name instead of the generic refArea label.
$ws.Range("C2").Value = "iaest-measure:residencia-comarca-nombre"
$ws.Range("D2").Value = "iaest-measure:nacionalidad-continente-nombre"
$ws.Range("F2").Value = "iaest-measure:residencia-ccaa-nombre"
$ws.Range("H2").Value = "iaest-measure:tipo-de-nacionalidad"
$ws.Range("K2").Value = "iaest-measure:sexo"
$ws.Range("M2").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("N2").Value = "iaest-measure:nacionalidad-area-nombre"

# Row 3: "dim" -> "medida" for the same set of curated columns.
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "medida"
$ws.Range("F3").Value = "medida"
$ws.Range("H3").Value = "medida"
$ws.Range("K3").Value = "medida"
$ws.Range("M3").Value = "medida"
$ws.Range("N3").Value = "medida"

# Row 4: datatype column; the two URI-based datatypes are replaced, and the rest of the
# skos:Concept columns (for the same curated set) become xsd:int as they're now measures.
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("K4").Value = "xsd:int"
$ws.Range("M4").Value = "xsd:int"
$ws.Range("N4").Value = "xsd:int"

# Row 5: mapping file references; only the "ano" mapping remains, the curated dimensions
# no longer use an external mapping workbook. Use Clear() so the cells are removed
# entirely rather than left behind as empty styled cells.
$ws.Range("D5").Clear()
$ws.Range("F5").Clear()
$ws.Range("H5").Clear()
$ws.Range("K5").Clear()
$ws.Range("N5").Clear()
